$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab03")

# 1. Insert a new column before the old AO column (shifts old AO -> AP, preserving its format/value)
$ws.Columns("AO:AO").Insert()

# 2. Update the report title (C1) and the "Average annual projected growth" header (now AP2)
$ws.Range("C1").Value = "Table 3: Annual population growth rate, 1990-2028"
$ws.Range("AP2").Value = "Average annual projected growth, 2023-28"

# 3. Set the new year header in AO2
$ws.Range("AO2").Value = 2028

# 4. Fill in the new 2028 growth-rate figures (AO) and the recalculated 2023-28 averages (AP)
# for each country/aggregate data row (rows 3 to 98).
$ws.Range("AO3").Value = 2.89259070649865
$ws.Range("AP3").Value = 2.9717274785394001
$ws.Range("AO4").Value = 1.46501282715565
$ws.Range("AP4").Value = 1.55813392781803
$ws.Range("AO5").Value = 1.1026893314727999
$ws.Range("AP5").Value = 1.04126686049189
$ws.Range("AO6").Value = 0.99164033701951004
$ws.Range("AP6").Value = 1.04644391076207
$ws.Range("AO7").Value = 2.5123843261929402
$ws.Range("AP7").Value = 2.5658697345564399
$ws.Range("AO8").Value = 2.69195808192215
$ws.Range("AP8").Value = 2.7647041212568699
$ws.Range("AO9").Value = 1.5820871457938701
$ws.Range("AP9").Value = 1.6236875231584
$ws.Range("AO10").Value = 0.94659291332992002
$ws.Range("AP10").Value = 1.0058088885441601
$ws.Range("AO11").Value = 2.60193042516324
$ws.Range("AP11").Value = 2.67253458666625
$ws.Range("AO12").Value = 1.9943621069375499
$ws.Range("AP12").Value = 2.0567872786027501
$ws.Range("AO13").Value = 2.0743143349136099
$ws.Range("AP13").Value = 2.1270602197931399
$ws.Range("AO14").Value = 2.5259639943046501
$ws.Range("AP14").Value = 2.5893985687769399
$ws.Range("AO15").Value = 2.4619236794633399
$ws.Range("AP15").Value = 2.5345556044657398
$ws.Range("AO16").Value = 3.1202566427631
$ws.Range("AP16").Value = 3.0889374531711602
$ws.Range("AO17").Value = 2.9461734959618
$ws.Range("AP17").Value = 3.0270585669090999
$ws.Range("AO18").Value = 2.1906515405715599
$ws.Range("AP18").Value = 2.2174347151837002
$ws.Range("AO19").Value = 3.1862099993069402
$ws.Range("AP19").Value = 3.2384354190141198
$ws.Range("AO20").Value = 2.1711171255647499
$ws.Range("AP20").Value = 2.26345733651161
$ws.Range("AO21").Value = 1.8334500821044299
$ws.Range("AP21").Value = 1.9059793676694099
$ws.Range("AO22").Value = 1.9208267941360999
$ws.Range("AP22").Value = 1.9389824171407899
$ws.Range("AO23").Value = 2.9354403904974
$ws.Range("AP23").Value = 2.9893280141633798
$ws.Range("AO24").Value = 1.67025968199079
$ws.Range("AP24").Value = 1.7453108992042199
$ws.Range("AO25").Value = 1.3177619545607899
$ws.Range("AP25").Value = 1.36023826950213
$ws.Range("AO26").Value = 1.95808186726947
$ws.Range("AP26").Value = 1.9045328722233199
$ws.Range("AO27").Value = 2.3494879623416498
$ws.Range("AP27").Value = 2.4372159960398099
$ws.Range("AO28").Value = 1.9436340227792801
$ws.Range("AP28").Value = 1.97565134677762
$ws.Range("AO29").Value = 2.2902152104265601
$ws.Range("AP29").Value = 2.35088977342917
$ws.Range("AO30").Value = 0.03586467980428
$ws.Range("AP30").Value = 0.073951119308989996
$ws.Range("AO31").Value = 2.1301040301952701
$ws.Range("AP31").Value = 2.2010187010445699
$ws.Range("AO32").Value = 0.44671694885955998
$ws.Range("AP32").Value = 0.50463404892658004
$ws.Range("AO33").Value = 2.9669867964191101
$ws.Range("AP33").Value = 3.0372006158139402
$ws.Range("AO34").Value = 1.83974240661979
$ws.Range("AP34").Value = 1.7707433355601301
$ws.Range("AO35").Value = 2.4021914706760898
$ws.Range("AP35").Value = 2.4985976007434898
$ws.Range("AO36").Value = 2.7719137427934899
$ws.Range("AP36").Value = 2.8521707829067799
$ws.Range("AO37").Value = 2.6376090546270201
$ws.Range("AP37").Value = 2.6984225584235402
$ws.Range("AO38").Value = 2.3933941363000599
$ws.Range("AP38").Value = 2.4597556373820999
$ws.Range("AO39").Value = 1.18476190571171
$ws.Range("AP39").Value = 1.32133393907923
$ws.Range("AO40").Value = 1.47957071681177
$ws.Range("AP40").Value = 1.5306256597110099
$ws.Range("AO41").Value = 0.98933432227787999
$ws.Range("AP41").Value = 1.0410347515910501
$ws.Range("AO42").Value = 2.6057891619780902
$ws.Range("AP42").Value = 2.65590220719423
$ws.Range("AO43").Value = 0.84266425318880001
$ws.Range("AP43").Value = 0.90845302232183001
$ws.Range("AO44").Value = 0.67503780768782995
$ws.Range("AP44").Value = 0.76168787269992
$ws.Range("AO45").Value = 1.27769017536452
$ws.Range("AP45").Value = 1.3487234351891699
$ws.Range("AO46").Value = 2.5480531721557198
$ws.Range("AP46").Value = 2.6172163415424401
$ws.Range("AO47").Value = 2.4049334415088701
$ws.Range("AP47").Value = 2.4694837632392099
$ws.Range("AO48").Value = 0.92518127247809001
$ws.Range("AP48").Value = 0.94366674704970999
$ws.Range("AO49").Value = 2.42718525839161
$ws.Range("AP49").Value = 2.4768911793682902
$ws.Range("AO50").Value = 2.3064534946980002
$ws.Range("AP50").Value = 2.3959028132729498
$ws.Range("AO51").Value = 1.8157645145624499
$ws.Range("AP51").Value = 1.8706177260715999
$ws.Range("AO52").Value = 2.2461197937407902
$ws.Range("AP52").Value = 2.3245743528273999
$ws.Range("AO53").Value = 2.0565240657253399
$ws.Range("AP53").Value = 2.1121914732523699
$ws.Range("AO54").Value = 2.1206337675749598
$ws.Range("AP54").Value = 2.1554025895855302
$ws.Range("AO55").Value = 3.00783701807732
$ws.Range("AP55").Value = 3.0608381959560398
$ws.Range("AO56").Value = 3.7384325482259899
$ws.Range("AP56").Value = 3.78346133193899
$ws.Range("AO57").Value = 2.2809699256261
$ws.Range("AP57").Value = 2.3364928256685
$ws.Range("AO58").Value = 2.4815394138618201
$ws.Range("AP58").Value = 2.53164788085083
$ws.Range("AO59").Value = 1.96897847476594
$ws.Range("AP59").Value = 2.0478169712932401
$ws.Range("AO60").Value = 2.19178807508038
$ws.Range("AP60").Value = 2.2382089719519098
$ws.Range("AO61").Value = 2.4003099771119798
$ws.Range("AP61").Value = 2.4527379385311798
$ws.Range("AO62").Value = 2.25831207275538
$ws.Range("AP62").Value = 2.3136583798521602
$ws.Range("AO63").Value = 0.51650227092372003
$ws.Range("AP63").Value = 0.55579199382750999
$ws.Range("AO64").Value = 0.66871683990358999
$ws.Range("AP64").Value = 0.71195074341992004
$ws.Range("AO65").Value = 0.61875852367416995
$ws.Range("AP65").Value = 0.66514680606855003
$ws.Range("AO66").Value = 0.85259978474318998
$ws.Range("AP66").Value = 0.88561014986778996
$ws.Range("AO67").Value = 2.2947830035001
$ws.Range("AP67").Value = 2.3509876362642501
$ws.Range("AO68").Value = 2.1697183188398301
$ws.Range("AP68").Value = 2.2205366146742298
$ws.Range("AO69").Value = 2.42064259351891
$ws.Range("AP69").Value = 2.47142957691617
$ws.Range("AO70").Value = 2.8805032655275302
$ws.Range("AP70").Value = 2.9387741541088399
$ws.Range("AO71").Value = 2.4003099771119798
$ws.Range("AP71").Value = 2.4527379385311798
$ws.Range("AO72").Value = 2.3435198651119702
$ws.Range("AP72").Value = 2.4085590237630901
$ws.Range("AO73").Value = 2.4940797835606401
$ws.Range("AP73").Value = 2.54708711239262
$ws.Range("AO74").Value = 1.06326937589478
$ws.Range("AP74").Value = 1.15686079420609
$ws.Range("AO75").Value = 2.7588067056466801
$ws.Range("AP75").Value = 2.8325664030873798
$ws.Range("AO76").Value = 0.74018451683870001
$ws.Range("AP76").Value = 0.79451319868351
$ws.Range("AO77").Value = 0.63254274081157003
$ws.Range("AP77").Value = 0.6745828236268
$ws.Range("AO78").Value = -0.10262777469910001
$ws.Range("AP78").Value = -0.19155038514799999
$ws.Range("AO79").Value = 0.23692414191041
$ws.Range("AP79").Value = 0.22519654161757
$ws.Range("AO80").Value = 2.21739778619146
$ws.Range("AP80").Value = 2.2821391218807299
$ws.Range("AO81").Value = 0.80204700189375
$ws.Range("AP81").Value = 0.84865977553409999
$ws.Range("AO82").Value = 2.2710643535862198
$ws.Range("AP82").Value = 2.3234891505905302
$ws.Range("AO83").Value = 0.49205738897327
$ws.Range("AP83").Value = 0.53087964327440995
$ws.Range("AO84").Value = 2.65145809237355
$ws.Range("AP84").Value = 2.71274047295667
$ws.Range("AO85").Value = 2.0789225361083501
$ws.Range("AP85").Value = 2.2744240031302501
$ws.Range("AO86").Value = 2.0583124590940298
$ws.Range("AP86").Value = 2.1147847197332799
$ws.Range("AO87").Value = 0.90669878612461996
$ws.Range("AP87").Value = 0.98192222993808997
$ws.Range("AO88").Value = 1.03233676465713
$ws.Range("AP88").Value = 1.0904939256039901
$ws.Range("AO89").Value = 0.13053906197331
$ws.Range("AP89").Value = 0.15559686847692999
$ws.Range("AO90").Value = 0.20074061059553999
$ws.Range("AP90").Value = 0.18669523602957
$ws.Range("AO91").Value = 2.6584443642795299
$ws.Range("AP91").Value = 2.7219061196761301
$ws.Range("AO92").Value = 1.1762423390891099
$ws.Range("AP92").Value = 1.2371394068225301
$ws.Range("AO93").Value = 1.3518798981957201
$ws.Range("AP93").Value = 1.38673439375814
$ws.Range("AO94").Value = 0.76161499185134995
$ws.Range("AP94").Value = 0.80236279698035995
$ws.Range("AO95").Value = 2.5518784483637802
$ws.Range("AP95").Value = 2.61249283401499
$ws.Range("AO96").Value = 1.3556487234738901
$ws.Range("AP96").Value = 1.41343969054546
$ws.Range("AO97").Value = 2.5228856454156499
$ws.Range("AP97").Value = 2.5816002394782398
$ws.Range("AO98").Value = 1.42627200968957
$ws.Range("AP98").Value = 1.5080500014995
